# Final Commit 4th Oct,2018
#
# The sheet tracks students' userId/userName/batch/password/academicyear/
# sem/dept/degree/class. This edit fills in the "academicyear" value for
# row 2 (column E, under the "academicyear" header) and leaves the active
# cell on the next column (F2), matching the author's recorded selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the academic-year value for the second data row.
$ws.Range("E2").Value = "2018-19_ODD"

# Move the active selection to F2 (was G2 before the edit).
$ws.Range("F2").Select()
